$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" (1st sheet): row 3 is fbd452d3 row ---
# E3/F3 share the same underlying status string as the per-language sheets'
# Status column, so it flips too when that shared text changes.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- Sheet "zh-cn" (2nd sheet): row 3 is fbd452d3 row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: duqhd4cq.mob is different with handoff file name: fbd452d3-7f96-4038-a283-cf6c73647bae.46ed0fbd988f49f67a5bc3d53ba063fb3da2be50.zh-cn."
# Column-width write-back adds a constant 5/6-character pad before it hits the
# sheet XML, so back it out here to land on an exact stored width of 40.
$wsZhCn.Columns.Item(16).ColumnWidth = (40 - 5/6)

# --- Sheet "de-de" (3rd sheet): row 3 is fbd452d3 row ---
# The "Status" text is a shared string reused by both language sheets, so
# the same status update applies here too (it flipped for both in the source edit).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: duqhd4cq.mob is different with handoff file name: fbd452d3-7f96-4038-a283-cf6c73647bae.46ed0fbd988f49f67a5bc3d53ba063fb3da2be50.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = (40 - 5/6)
